# Updated cryptos list on Sat Jan  6 22:29:31 UTC 2024 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) columns on the
# cryptos table to the latest scraped snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new Price text (column D). $null means "leave Price unchanged".
$prices = @{
    2  = "44.266.05"
    3  = "2.238.59"
    4  = "1.01"
    5  = "307.38"
    6  = "94.46"
    7  = "0.571"
    8  = "1.01"
    9  = "0.523"
    10 = "34.47"
    11 = "0.0809"
    12 = "7.16"
    14 = "2.579.90"
    15 = "2.317.09"
    16 = "0.826"
    17 = "13.45"
    18 = "44.003.32"
    19 = "0.0₃0963"
    20 = "6.36"
    21 = "12.10"
    22 = "65.50"
    23 = "236.48"
    24 = "2.95"
    27 = "39.75"
    29 = "9.86"
    31 = "5.85"
    32 = "153.55"
    33 = "0.0792"
    35 = "3.11"
    37 = "0.108"
    38 = "1.74"
    39 = "3.48"
    41 = "14.25"
    44 = "1.733.54"
    45 = "82.52"
    46 = "0.191"
    47 = "99.16"
    48 = "4.91"
    49 = "1.60"
    50 = "8.06"
    51 = "54.60"
}

# row -> new Volume(1h) percentage (column E, without the padding spaces).
# $null means "leave Volume(1h) unchanged".
$volumes = @{
    2  = "+0.33%"
    3  = "-0.26%"
    4  = "+0.25%"
    5  = "-2.70%"
    6  = "-5.02%"
    7  = "-0.64%"
    8  = "+0.36%"
    9  = "-2.19%"
    10 = "-5.01%"
    11 = "-1.88%"
    12 = "-3.15%"
    13 = "-0.04%"
    14 = "-0.28%"
    15 = "+3.15%"
    16 = "-2.26%"
    17 = "-3.91%"
    18 = "+0.08%"
    19 = "-1.90%"
    21 = "-8.06%"
    22 = "-0.08%"
    23 = "-0.80%"
    24 = "-1.53%"
    25 = "-1.10%"
    26 = "-0.04%"
    27 = "+7.31%"
    28 = "+4.05%"
    29 = "-2.67%"
    30 = "-0.12%"
    31 = "-2.47%"
    32 = "-1.39%"
    33 = "-5.66%"
    34 = "-2.37%"
    35 = "-6.32%"
    36 = "+1.53%"
    37 = "-2.11%"
    38 = "-8.06%"
    39 = "-1.72%"
    40 = "-4.77%"
    41 = "-7.03%"
    42 = "-3.63%"
    43 = "+0.29%"
    44 = "+1.84%"
    45 = "-1.29%"
    46 = "-2.40%"
    47 = "-2.69%"
    48 = "-5.51%"
    49 = "-0.99%"
    50 = "-0.65%"
    51 = "-3.60%"
}

foreach ($row in $prices.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $text = $prices[$row]

    # The Price column stores plain text (e.g. "1.01", "44.266.05"). Values
    # that parse as a plain number would otherwise be auto-coerced to a
    # numeric cell by `.Value`, so force the cell to Text first whenever the
    # new price looks like a bare number ("1.01"), same as it would need to
    # in real Excel. Values containing thousands separators (two dots) or
    # other non-numeric characters don't need this.
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $text
}

foreach ($row in $volumes.Keys) {
    $ws.Cells.Item($row, 5).Value = "  " + $volumes[$row] + "  "
}
